$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(876).Insert()

# Force column A's new value to be stored as text (it looks like a date and
# would otherwise get auto-converted to a date serial number), then restore
# the default "Normal" style so no stray number-format style is left behind.
$ws.Cells.Item(876, 1).NumberFormat = "@"
$ws.Cells.Item(876, 1).Value = "2026/02/28"
$ws.Cells.Item(876, 1).Style = "Normal"

$ws.Cells.Item(876, 2).Value = "土"
$ws.Cells.Item(876, 3).Value = 13
$ws.Cells.Item(876, 4).Value = 31
